# Refine row number dependant element location to parameterized settings
#
# The "Location" worksheet lists Key/Type/Value rows describing how the
# TestOnTankLibrary test-suite finds elements on the page. The XPath
# expressions that pointed at a specific (hard-coded) table row - e.g.
# "Home.List.All.Data2.*" / "Home.List.All.Data3.*" using literal
# tr[2] / tr[3] indices - are collapsed into a single, reusable,
# parameterized "Home.List.All.Data.*" location keyed on tr[{0}], and two
# new locations (row Edit/Delete buttons) plus an "Edit.Name" key group are
# appended at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Location")

# --- Rewrite the existing parameterized-row block (rows 9-14) ------------

$ws.Range("A9").Value  = "Home.List.All.Data.Name"
$ws.Range("B9").Value  = "XPath"
$ws.Range("C9").Value  = "//*[@id='tanklist']/tbody/tr[{0}]/td[2]"

$ws.Range("A10").Value = "Home.List.All.Data.Stage"
$ws.Range("B10").Value = "XPath"
$ws.Range("C10").Value = "//*[@id='tanklist']/tbody/tr[{0}]/td[3]"

$ws.Range("A11").Value = "Home.List.All.Data.Type"
$ws.Range("B11").Value = "XPath"
$ws.Range("C11").Value = "//*[@id='tanklist']/tbody/tr[{0}]/td[4]"

$ws.Range("A12").Value = "Home.Desc.Name"
$ws.Range("B12").Value = "XPath"
$ws.Range("C12").Value = "//*[@id='imganddesc_div']/div/h4"

$ws.Range("A13").Value = "Detail.Name"
$ws.Range("B13").Value = "XPath"
$ws.Range("C13").Value = "//dt[text()='Name']/following-sibling::dd[1]"

$ws.Range("A14").Value = "Home.List.All.Data.Btn.Edit"
$ws.Range("B14").Value = "XPath"
$ws.Range("C14").Value = "//*[@id='tanklist']/tbody/tr[{0}]/td[1]//i[@class='fa fa-edit']"

# --- Append the two new rows, copying the existing striped row format ----

$ws.Range("A13:G13").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Rows.Item(15).RowHeight = 20.25

$ws.Range("A14:G14").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Rows.Item(16).RowHeight = 20.25

$ws.Range("A15").Value = "Home.List.All.Data.Btn.Delete"
$ws.Range("B15").Value = "XPath"
$ws.Range("C15").Value = "//*[@id='tanklist']/tbody/tr[{0}]/td[1]//i[@class='fa fa-trash']"

$ws.Range("A16").Value = "Edit.Name"
$ws.Range("B16").Value = "Id"
$ws.Range("C16").Value = "Name"

$ws.Range("D16").Value = ""
$ws.Range("E16").Value = ""
$ws.Range("F16").Value = ""
$ws.Range("G16").Value = ""
